# MODELO_PROPOSTA.docx edit:
#  1. Move the stray "_GoBack" bookmark away from the CNPJ run (inside the
#     Artista qualification paragraph) to the blank paragraph that sits
#     between the "OBSERVACOES" table and the "VALIDADE DA PROPOSTA" table.
#     (Word re-targets a bookmark on Add() if the name already exists, so a
#     single Bookmarks.Add with the new location both removes the old
#     placement and creates the new one.)
#  2. Resize the two-column key/value tables (FORMA DE PAGAMENTO,
#     OBSERVACOES, VALIDADE DA PROPOSTA) so both columns share a common
#     width (2240 / 6248 twips == 112 / 312.4 points).

$d = $word.ActiveDocument

# --- 1. Relocate the "_GoBack" bookmark -------------------------------

$endRange = $d.Content
$endRange.Find.Execute("jurídico.", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null

$startRange = $d.Content
$startRange.Find.Execute("VALIDADE DA PROPOSTA:", $false, $false, $false, `
    $false, $false, $true, 1, $false, "", 0) | Out-Null

$gap = $d.Range($endRange.End, $startRange.Start)
$d.Bookmarks.Add("_GoBack", $gap) | Out-Null

# --- 2. Resize the key/value tables ------------------------------------

function Set-TwoColumnWidths($table, $col1Twips, $col2Twips) {
    $table.Columns.Item(1).Width = $col1Twips / 20.0
    $table.Columns.Item(2).Width = $col2Twips / 20.0
}

foreach ($i in 1..$d.Tables.Count) {
    $t = $d.Tables.Item($i)
    $label = $t.Cell(1, 1).Range.Text
    if ($label -like "FORMA DE PAGAMENTO:*" -or `
        $label -like "OBSERVA*ES:*" -or `
        $label -like "VALIDADE DA PROPOSTA:*") {
        Set-TwoColumnWidths $t 2240 6248
    }
}
